$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Week number
$ws.Range("E1").Value = 11

# Task rows - "Stage" column filled first for all rows
$ws.Range("A3").Value = "Project Build"
$ws.Range("A4").Value = "Project Build"
$ws.Range("A5").Value = "Project Build"
$ws.Range("A6").Value = "Project Build"

# Totals row label
$ws.Range("A14").Value = "Cumulative Total: 220"
$ws.Range("D14").Value = 20

# Task descriptions
$ws.Range("B3").Value = "Implementation of final feature test"
$ws.Range("B4").Value = "Unit Testing"

# Name
$ws.Range("C1").Value = "Richard Dobson"

$ws.Range("B5").Value = "Integration Testing between front end and back end prog"
$ws.Range("B6").Value = "Finish Documentation"

# Numeric hour columns
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5

# Selection moves to B3
$ws.Range("B3").Select()
